$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Critical issues")

$ws.Range("A4").Value = "Excluded by ‘noindex’ tag"
$ws.Range("B4").Value = "Website"
$ws.Range("C4").Value = "Not Started"
$ws.Range("D4").Value = 15.0

$ws.Range("A5").Value = "Blocked by robots.txt"
$ws.Range("B5").Value = "Website"
$ws.Range("C5").Value = "Not Started"
$ws.Range("D5").Value = 1.0

$ws.Range("A6").Value = "Page with redirect"
$ws.Range("B6").Value = "Website"
$ws.Range("C6").Value = "Started"
$ws.Range("D6").Value = 2.0

$ws.Range("A7").Value = "Page with redirect"
$ws.Range("B7").Value = "Website"
$ws.Range("C7").Value = "Started"
$ws.Range("D7").Value = 1.0

$ws.Range("A8").Value = "Duplicate, Google chose different canonical than user"
$ws.Range("B8").Value = "Google systems"
$ws.Range("C8").Value = "Started"
$ws.Range("D8").Value = 36.0
